# Update Jordan Leggett TE 2018 sheet: insert "height" and "weight" columns
# between the existing "fumbles" and "fantasy points" columns, shifting the
# "fantasy points" data (previously column E) out to column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing "fantasy points" values, currently stored in column E,
# for each data row (row number -> value).
$fantasyPoints = @{
    2  = 0
    3  = 0
    4  = 7.3
    5  = 1
    6  = 0
    7  = 0
    8  = 0
    9  = 0.6
    10 = 0.6
    11 = 2.4
    12 = 1.8
    13 = 1.3
    14 = 0.8
    15 = 0.9
    16 = 0.7
}

# Move the "fantasy points" column from E to G (header + values), copying
# the existing header's formatting (bold font + border) so the new header
# cell matches the rest of the header row.
$ws.Range("D1").Copy() | Out-Null
$ws.Range("G1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("G1").Value = "fantasy points"
foreach ($r in $fantasyPoints.Keys) {
    $ws.Cells.Item($r, 7).Value = $fantasyPoints[$r]
}

# New headers for the inserted columns. E1 already carries the header
# style from the original "fantasy points" header; F1 needs it copied
# over since it is a brand-new cell.
$ws.Range("D1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("E1").Value = "height"
$ws.Range("F1").Value = "weight"

# Fill the new height/weight columns with the scraped values for every
# data row (rows 2-16).
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 5).Value = 6.416666666666667
    $ws.Cells.Item($r, 6).Value = 258
}
